# Updated cryptos list on Sat Mar 16 06:31:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to be treated as text so that values
# like "1.00" or "58.34" aren't silently coerced into numbers (which would
# drop meaningful trailing/format digits), matching the inlineStr cells in
# the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Row 2 : Bitcoin ---
$ws.Range("D2").Value = "69.043.84"
$ws.Range("E2").Value = "  +0.88%  "

# --- Row 3 : Ethereum ---
$ws.Range("D3").Value = "3.720.50"
$ws.Range("E3").Value = "  -0.55%  "

# --- Row 4 : TetherUSD ---
$ws.Range("E4").Value = "  +0.42%  "

# --- Row 5 : BNB ---
$ws.Range("D5").Value = "611.37"
$ws.Range("E5").Value = "  +4.93%  "

# --- Row 6 : Solana ---
$ws.Range("D6").Value = "190.35"
$ws.Range("E6").Value = "  +6.19%  "

# --- Row 7 ---
$ws.Range("E7").Value = "  +0.02%  "

# --- Row 8 ---
$ws.Range("E8").Value = "  +0.16%  "

# --- Row 9 ---
$ws.Range("E9").Value = "  -0.67%  "

# --- Row 10 ---
$ws.Range("E10").Value = "  -4.75%  "

# --- Row 11 ---
$ws.Range("D11").Value = "58.34"
$ws.Range("E11").Value = "  +8.15%  "

# --- Row 12 ---
$ws.Range("E12").Value = "  -4.67%  "

# --- Row 13 ---
$ws.Range("E13").Value = "  -2.58%  "

# --- Row 14 ---
$ws.Range("D14").Value = "4.313.47"
$ws.Range("E14").Value = "  +0.42%  "

# --- Row 15 ---
$ws.Range("D15").Value = "3.720.86"
$ws.Range("E15").Value = "  -0.01%  "

# --- Row 16 ---
$ws.Range("D16").Value = "19.32"
$ws.Range("E16").Value = "  -1.30%  "

# --- Row 17 ---
$ws.Range("E17").Value = "  -0.34%  "

# --- Row 18 ---
$ws.Range("E18").Value = "  -1.06%  "

# --- Row 19 ---
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  -1.51%  "

# --- Row 20 ---
$ws.Range("D20").Value = "68.877.88"
$ws.Range("E20").Value = "  +1.29%  "

# --- Row 21 ---
$ws.Range("D21").Value = "411.55"
$ws.Range("E21").Value = "  -0.49%  "

# --- Row 22 ---
$ws.Range("E22").Value = "  -0.30%  "

# --- Row 23 ---
$ws.Range("D23").Value = "89.43"
$ws.Range("E23").Value = "  +0.11%  "

# --- Row 24 ---
$ws.Range("E24").Value = "  -2.08%  "

# --- Row 25 ---
$ws.Range("E25").Value = "  -0.75%  "

# --- Row 26 ---
$ws.Range("D26").Value = "10.90"
$ws.Range("E26").Value = "  +0.19%  "

# --- Row 27 ---
$ws.Range("D27").Value = "6.05"
$ws.Range("E27").Value = "  +0.93%  "

# --- Row 28 ---
$ws.Range("D28").Value = "3.80"
$ws.Range("E28").Value = "  -1.65%  "

# --- Row 29 ---
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  -0.09%  "

# --- Row 30 ---
$ws.Range("D30").Value = "33.13"
$ws.Range("E30").Value = "  -0.62%  "

# --- Row 31 ---
$ws.Range("D31").Value = "7.48"
$ws.Range("E31").Value = "  -9.21%  "

# --- Row 32 ---
$ws.Range("D32").Value = "12.71"
$ws.Range("E32").Value = "  -0.55%  "

# --- Row 33 ---
$ws.Range("E33").Value = "  +3.37%  "

# --- Row 34 ---
$ws.Range("D34").Value = "45.90"
$ws.Range("E34").Value = "  +2.76%  "

# --- Row 35 ---
$ws.Range("D35").Value = "625.69"
$ws.Range("E35").Value = "  +1.75%  "

# --- Row 36 ---
$ws.Range("D36").Value = "65.62"
$ws.Range("E36").Value = "  -0.15%  "

# --- Row 37 ---
$ws.Range("D37").Value = "0.414"
$ws.Range("E37").Value = "  +1.76%  "

# --- Row 38 ---
$ws.Range("D38").Value = "0.0₃0819"
$ws.Range("E38").Value = "  -13.18%  "

# --- Row 39 ---
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.14%  "

# --- Row 40 ---
$ws.Range("E40").Value = "  +0.55%  "

# --- Row 41 ---
$ws.Range("E41").Value = "  +1.65%  "

# --- Row 42 ---
$ws.Range("E42").Value = "  -2.08%  "

# --- Row 43 ---
$ws.Range("E43").Value = "  +0.19%  "

# --- Row 44 ---
$ws.Range("D44").Value = "2.62"
$ws.Range("E44").Value = "  -1.70%  "

# --- Row 45 ---
$ws.Range("E45").Value = "  +3.01%  "

# --- Row 46 : Maker ---
$ws.Range("D46").Value = "2.855.44"
$ws.Range("E46").Value = "  +4.02%  "

# --- Row 47 : WEMIXToken ---
$ws.Range("E47").Value = "  -0.12%  "

# --- Row 48 : THORChain ---
$ws.Range("E48").Value = "  -5.54%  "

# --- Row 49 : was ApeXProtocol, now Monero (rows 49/50 swapped) ---
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "143.70"
$ws.Range("E49").Value = "  -0.18%  "

# --- Row 50 : was Monero, now ApeXProtocol ---
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "3.10"
$ws.Range("E50").Value = "  -2.01%  "

# --- Row 51 : dogwifhat ---
$ws.Range("D51").Value = "2.58"
$ws.Range("E51").Value = "  -21.40%  "
